$d = $word.ActiveDocument

# --- Edit 1: extend the first paragraph with a colored annotation ---
$p1 = $d.Paragraphs(1)
$r = $p1.Range

# Two trailing spaces appended to the existing (uncolored) run.
$r.InsertAfter("  ")

# "(This is a change " + en-dash + " Ve" in red (C00000)
$t1 = "(This is a change " + [char]0x2013 + " Ve"
$r.InsertAfter($t1)
$run1End = $r.End - 1
$run1Start = $run1End - $t1.Length
$cr1 = $d.Range($run1Start, $run1End)
$cr1.Font.Color = 192

# "rsion for branch alternate" in red (C00000)
$t2 = "rsion for branch alternate"
$r.InsertAfter($t2)
$run2End = $r.End - 1
$run2Start = $run2End - $t2.Length
$cr2 = $d.Range($run2Start, $run2End)
$cr2.Font.Color = 192

# ")" in red (C00000)
$t3 = ")"
$r.InsertAfter($t3)
$run3End = $r.End - 1
$run3Start = $run3End - $t3.Length
$cr3 = $d.Range($run3Start, $run3End)
$cr3.Font.Color = 192

# --- Edit 2: append a new empty, shaded paragraph after the final paragraph ---
$d.Content.Find.Execute("we are free at last.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "we are free at last.^p", 2)

$newp = $d.Paragraphs($d.Paragraphs.Count)
$newp.Style = -1

$sh = $newp.Shading
$sh.Texture = 0
$sh.ForegroundPatternColor = -16777216
$sh.BackgroundPatternColor = 16382457
